$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

$ws.Range("I7").Value = "nan"
$ws.Range("I8").Value = "✅"
